$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.360.33"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").Value = "3.332.94"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'553.72"
$ws.Range("E5").Value = "  +0.71%  "

$ws.Range("D6").Value = "'173.74"
$ws.Range("E6").Value = "  +0.73%  "

$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = "  +2.08%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "3.322.37"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("E10").Value = "  +6.92%  "

$ws.Range("D11").Value = "'0.634"
$ws.Range("E11").Value = "  +2.32%  "

$ws.Range("D12").Value = "'53.42"
$ws.Range("E12").Value = "  +0.60%  "

$ws.Range("D13").Value = "'0.0000279"
$ws.Range("E13").Value = "  +2.36%  "

$ws.Range("E14").Value = "  +1.24%  "

$ws.Range("D15").Value = "3.861.36"
$ws.Range("E15").Value = "  +0.04%  "

$ws.Range("E16").Value = "  +3.20%  "

$ws.Range("D17").Value = "'18.17"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").Value = "3.343.25"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'11.77"
$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "64.373.79"
$ws.Range("E20").Value = "  +0.70%  "

$ws.Range("D21").Value = "'0.989"
$ws.Range("E21").Value = "  +1.94%  "

$ws.Range("D22").Value = "'453.22"
$ws.Range("E22").Value = "  +6.83%  "

$ws.Range("D23").Value = "'5.00"
$ws.Range("E23").Value = "  +6.96%  "

$ws.Range("D24").Value = "'4.07"
$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").Value = "'87.96"
$ws.Range("E25").Value = "  +4.87%  "

$ws.Range("E26").Value = "  +5.16%  "

$ws.Range("D27").Value = "'2.89"
$ws.Range("E27").Value = "  +3.00%  "

$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("D29").Value = "'31.07"
$ws.Range("E29").Value = "  +4.92%  "

$ws.Range("D30").Value = "'8.61"
$ws.Range("E30").Value = "  +0.48%  "

$ws.Range("D31").Value = "'6.53"
$ws.Range("E31").Value = "  -2.40%  "

$ws.Range("D32").Value = "'11.43"
$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("D33").Value = "'62.00"
$ws.Range("E33").Value = "  +6.56%  "

$ws.Range("D34").Value = "'569.21"
$ws.Range("E34").Value = "  -4.16%  "

$ws.Range("E35").Value = "  +0.55%  "

$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("E37").Value = "  -0.24%  "

$ws.Range("E38").Value = "  +1.01%  "

$ws.Range("D39").Value = "'35.45"
$ws.Range("E39").Value = "  +0.43%  "

$ws.Range("D40").Value = "'0.368"
$ws.Range("E40").Value = "  +1.07%  "

$ws.Range("D41").Value = "0.0₃0731"
$ws.Range("E41").Value = "  -2.28%  "

$ws.Range("D42").Value = "3.069.10"
$ws.Range("E42").Value = "  -0.99%  "

$ws.Range("E43").Value = "  +2.83%  "

$ws.Range("E44").Value = "  -0.60%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.20"
$ws.Range("E45").Value = "  +0.54%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.46"
$ws.Range("E46").Value = "  +0.83%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.134"
$ws.Range("E47").Value = "  +4.25%  "

$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "  -0.09%  "

$ws.Range("D49").Value = "'141.07"
$ws.Range("E49").Value = "  +6.43%  "

$ws.Range("D50").Value = "'2.50"
$ws.Range("E50").Value = "  -2.65%  "

$ws.Range("D51").Value = "'8.16"
$ws.Range("E51").Value = "  +0.20%  "
